$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs "2." / "9" / "+" (identical formatting)
# into a single run "2.9+". Find/Replace across the whole matched range
# naturally collapses it into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2.9+", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2.9+", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: update the Stackoverflow profile URL text. Only the trailing
# "parthiban-s" -> "parthiban-soundram" portion is searched/replaced (rather
# than the whole URL) so the match stays strictly inside the existing run
# and its InternetLink run style/formatting is preserved.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "parthiban-s", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "parthiban-soundram", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: rework the last bullet (previously "Github <link>") into a
# "Portfolio  https://parthibansoundram.github.io" line, with a bookmark
# around "Portfolio", Internet-Link styled/no-underline runs, and an extra
# blank paragraph added right after it.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Github*") {
        $target = $p
        break
    }
}

# NOTE: Word's InsertXML does not round-trip a <w:rStyle> character-style
# reference (it silently drops it, while keeping direct formatting such as
# <w:color>/<w:u>). So the XML below only carries the *direct* "color"
# override; the InternetLink run style and the "no underline" direct
# override are applied afterwards with Range.Style / Range.Font.Underline.
$xml = '<w:p ' + `
    'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
    'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
    '<w:pPr>' + `
      '<w:pStyle w:val="Normal1"/>' + `
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>' + `
      '<w:spacing w:lineRule="auto" w:line="360"/>' + `
      '<w:ind w:left="840" w:hanging="420"/>' + `
      '<w:jc w:val="both"/>' + `
      '<w:rPr></w:rPr>' + `
    '</w:pPr>' + `
    '<w:bookmarkStart w:id="0" w:name="__DdeLink__262_4104747235"/>' + `
    '<w:r>' + `
      '<w:rPr><w:color w:val="000000"/></w:rPr>' + `
      '<w:t>Portfolio</w:t>' + `
    '</w:r>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r>' + `
      '<w:rPr><w:color w:val="000000"/></w:rPr>' + `
      '<w:t xml:space="preserve"> </w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
      '<w:rPr></w:rPr>' + `
      '<w:t xml:space="preserve">  </w:t>' + `
    '</w:r>' + `
    '<w:hyperlink r:id="rId6">' + `
      '<w:r>' + `
        '<w:rPr></w:rPr>' + `
        '<w:t>https://</w:t>' + `
      '</w:r>' + `
      '<w:r>' + `
        '<w:rPr></w:rPr>' + `
        '<w:t>parthibansoundram.github.io</w:t>' + `
      '</w:r>' + `
    '</w:hyperlink>' + `
    '</w:p>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:pPr>' + `
      '<w:pStyle w:val="Normal1"/>' + `
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' + `
      '<w:spacing w:lineRule="auto" w:line="360"/>' + `
      '<w:ind w:left="840" w:hanging="0"/>' + `
      '<w:jc w:val="both"/>' + `
      '<w:rPr></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r><w:rPr></w:rPr></w:r>' + `
    '</w:p>'

if ($target -ne $null) {
    $target.Range.InsertXML($xml)

    # Re-locate the freshly-inserted "Portfolio ... github.io" paragraph and
    # restyle its text (everything except the trailing paragraph mark) as
    # InternetLink with no underline, matching the target formatting.
    $d2 = $word.ActiveDocument
    $newPara = $null
    for ($j = 1; $j -le $d2.Paragraphs.Count; $j++) {
        $pp = $d2.Paragraphs.Item($j)
        if ($pp.Range.Text -like "Portfolio*") {
            $newPara = $pp
            break
        }
    }
    if ($newPara -ne $null) {
        $full = $newPara.Range
        $textRange = $d2.Range($full.Start, $full.End - 1)
        $textRange.CharacterStyle = "InternetLink"
        $textRange.Font.Underline = 0
    }
}
